$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    ,@('UNIDEB_31', 'H', 'Powertrain System RD Team 2, Sprint 1', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_32', 'I', 'I will use the abbreviation PS for Powertrain System', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_33', 'R', 'When the PS receive a Gear position signal (id 7) the value shoud be saved (if the following conditions meet) and from then it shall be deemed that the car is in that gear.', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_34', 'R', 'If the car''s speed is 0 km/h then the saveable gear positions are: P, R, N, D (for the meaning of saveable see the requirement above)', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_35', 'R', 'If the car''s speed is bigger than 0 km/h forwards ( so the car is moving forward direction) then the only saveable gear positions are N and D', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_36', 'R', 'Every time if the received gear position signal''s value is equals to the currently saved one (eg the driver pressed D twice or more consequently) then the signal shall not take effect on the gear positions', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_37', 'R', 'If the car''s speed is bigger than 0 km/h backwards (so the car is moving bacward direction) then the only saveable gear positions are N and R', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_38', 'R', 'If the PS receive a Headlight signal (id  10) then it should save the received value.', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_39', 'R', 'The default value of the Headlight shall be false (so by default, the headlights are turned off)', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_40', 'R', 'The default value of the Gear position shall be P (so at the beginning of the application the car should be in P gear)', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_41', 'R', 'The default value of the Steering wheel angle shall be 0°', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_42', 'R', 'When the PS receive a Steering wheel angle signal (id 6) then the received value should be added to the saved one until it reches the minimum (-720°) or the maximum (720°) value. When the value after the addition exceeds the minimum value, the value should be the minimum. When the value after the addition exceeds the maximum value, the value shoud be the maximum.', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_43', 'R', 'The default value of the Indicator shall be 0 (NONE) at the beginning of the application', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_44', 'R', 'When the PS receive an Indicator signal (id 3) the logical value shoud be saved (if the following conditions meet)', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_45', 'R', 'Every time if the received Indicator signal''s value is equals to the currently saved one (eg the currently saved is 1 (RIGHT) and the received is also 1 (RIGHT) ) then the signal shall not take effect on the indicator''s saved value (nothing should happen)', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_46', 'R', 'If the PS receive an Indicator signal with the value of 0 (NONE) then both indicators shall be turned off (it should not blink)', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_47', 'R', 'If the PS receive an Indicator signal with the value of 1 (RIGHT) the right indicator should and only the right should start to blink (with 60Hz) and the left indicator shall be turned off.', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_48', 'R', 'If the PS receive an Indicator signal with the value of 2 (LEFT) the left indicator should and only the ledt should start to blink (with 60Hz) and the right indicator shall be turned off.', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_49', 'R', 'If the PS receive an Indicator signal  with the value of 3 (EMERGENCY) both left and right indicators should start to blink in synchronized (both turn on and off at the same time while blinking) (with 60Hz)', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_50', 'R', 'If the PS receive a Gas pedal position signal (id 7) and the car is in P gear (see above) then the received signal should not take any effect.', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_51', 'R', 'If the PS receive a Gas pedal position signal and the car is in N gear then the car motor''s rpm should be calculated based on the received signal''s value (which is in percentage) (eg if the received signal value is 50% then the rpm should be the 50% of the maximum rpm value) and this calculated value should be considered when the PS sends the motor rpm signal (id 9)', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_52', 'R', 'If the PS receive a Gas pedal position signal and the car is in R gear then the car motor''s rpm should be calculated based on the received signal''s value (which is in percentage) (eg if the received signal value is 50% then the rpm should be the 50% of the maximum rpm value) and this calculated value should be considered when the PS sends the motor rpm signal (id 9) and the following also has to be done:
Based on the saved steering wheel angle the car should start to accelerate with increasing the speed of the car to bacward direction and this should be done by calculating the car position x and y values, the car angle value and the car''s speed and sending these values by signals to the VirtualFunctionBus
The acceleration should continue while receiving the signal and when the signal is not received anymore, the acceleration should stop and the car should gradually slow down until it reaches the speed of 0 km/h.
The car''s speed should not exceed the maximum of 120 km/h', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_53', 'R', 'If the PS receive a Gas pedal position signal and the car is in D gear then the car motor''s rpm should be calculated based on the received signal''s value (which is in percentage) (eg if the received signal value is 50% then the rpm should be the 50% of the maximum rpm value) and this calculated value should be considered when the PS sends the motor rpm signal (id 9) and the following also has to be done:
Based on the saved steering wheel angle the car should start to accelerate with increasing the speed of the car to forward direction and this should be done by calculating the car position x and y values, the car angle value and the car''s speed and sending these values by signals to the VirtualFunctionBus
The acceleration should continue while receiving the signal and when the signal is not received anymore, the acceleration should stop and the car should gradually slow down until it reaches the speed of 0 km/h.
The car''s speed should not exceed the maximum of 120 km/h', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
    ,@('UNIDEB_54', 'R', 'If the PS receive a Brake pedal position signal (id 5) and the car''s speed is bigger than 0 km/h then the car''s speed should be decreased by the third of the current speed until it reaches the 0 km/h speed. If the car''s speed is 0 km/h then the brake pedal position signal should not take any effect.', 'New', 'https://trello.com/c/htzBlck2/6-rd-powertrain-system', 'Powertrain System', 'Team 2', 'Sprint 1')
)

$startRow = 33
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $srcRow = $r - 1
    $ws.Range("A" + $srcRow + ":H" + $srcRow).Copy()
    $ws.Range("A" + $r + ":H" + $r).PasteSpecial(-4122)
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

$excel.CutCopyMode = $false

$ws.Range("G33:G56").Select()

